# Add a second "email banner" slide to the deck.
#
# The new slide 2 is a duplicate of slide 1 (same picture, overlay
# rectangle and the two "Delta Delta Sigma" / "THE UNIVERSITY OF
# WASHINGTON'S PRE-DENTAL SOCIETY" text boxes), except the two text
# boxes get a slightly larger left inset (95.04pt == 1207008 EMU
# instead of 90pt == 1143000 EMU).

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Duplicate slide 1 -> becomes slide 2, inserted right after slide 1.
# This clones the picture + rectangle + both text boxes (shapes,
# geometry, formatting and text all copied verbatim), wires up a new
# slide part + relationships, and appends the new slide id to the
# presentation's slide list.
$s2 = $s1.Duplicate()
$slide2 = $p.Slides.Item(2)

# Bump the left text inset on the two text boxes only (msoTextBox ==
# 17; this excludes the picture and the semi-transparent overlay
# rectangle, which keep their original bodyPr/no lIns override) from
# 90pt to 95.04pt, matching the target markup's lIns="1207008".
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.Type -eq 17) {
        $shp.TextFrame.MarginLeft = 95.04
    }
}
